$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Baggy"
$ws.Range("B5").Value = "350 TL "
$ws.Range("C5").Value = "Jeans"
$ws.Range("D5").Value = "BAG1.jpg"
$ws.Range("E5").Value = "100% Pamuk"
$ws.Range("F5").Value = "Var"

$ws.Range("F5").Select()
